# WebForm User Assignment execution
# Re-assigns randomly generated phone numbers (column F, "PN_Value") for every
# participant row, and updates the Match2UserPos / UnMatchUserPos execution
# results (AN2 / AO2) produced by the webform run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New phone numbers (PN_Value) generated for this execution run.
$phoneNumbers = @{
    2  = "9840043653"
    3  = "9840021395"
    4  = "9840016496"
    5  = "9840052546"
    6  = "9840005196"
    7  = "9840078560"
    8  = "9840012905"
    9  = "9840038776"
    10 = "9840040363"
    11 = "9840090670"
    12 = "9840058219"
    13 = "9840009959"
    14 = "9840011059"
    15 = "9840049423"
    16 = "9840083865"
    17 = "9840031930"
    18 = "9840006741"
}

# Force the phone number column to be written as text (matching the
# existing "General" / text-stored-number formatting used in this sheet)
# rather than being auto-converted to a numeric value by Excel.
$phoneRange = $ws.Range("F2:F18")
$phoneRange.NumberFormat = "@"
foreach ($row in $phoneNumbers.Keys) {
    $ws.Cells.Item($row, 6).Value = $phoneNumbers[$row]
}
$phoneRange.ClearFormats()

# Update the webform execution summary counts.
$resultRange = $ws.Range("AN2:AO2")
$resultRange.NumberFormat = "@"
$ws.Range("AN2").Value = "1"
$ws.Range("AO2").Value = "0"
$resultRange.ClearFormats()
